$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Clcf1"
$ws.Range("C2").Value = "Il6st"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.679012
$ws.Range("H2").Value = 5.037036000000001
$ws.Range("I2").Value = 0.1178149724053671
$ws.Range("J2").Value = 0.1178149724053671
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 36.71344366666667
$ws.Range("N2").Value = 110.140331
$ws.Range("O2").Value = 0.2081992981130139
$ws.Range("P2").Value = 0.2081992981130138
$ws.Range("Q2").Value = 61.64231247765735
$ws.Range("R2").Value = 554.7808122989161
$ws.Range("S2").Value = 0.02452899456200153
$ws.Range("T2").Value = 0.02452899456200153

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Clcf1"
$ws.Range("C3").Value = "Il6st"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.679012
$ws.Range("H3").Value = 5.037036000000001
$ws.Range("I3").Value = 0.1178149724053671
$ws.Range("J3").Value = 0.1178149724053671
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 88.87708033333332
$ws.Range("N3").Value = 266.631241
$ws.Range("O3").Value = 0.5040155293450301
$ws.Range("P3").Value = 0.50401552934503
$ws.Range("Q3").Value = 149.2256844046307
$ws.Range("R3").Value = 1343.031159641676
$ws.Range("S3").Value = 0.05938057568166124
$ws.Range("T3").Value = 0.05938057568166123

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Clcf1"
$ws.Range("C4").Value = "Il6st"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.679012
$ws.Range("H4").Value = 5.037036000000001
$ws.Range("I4").Value = 0.1178149724053671
$ws.Range("J4").Value = 0.1178149724053671
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 14.22727166666667
$ws.Range("N4").Value = 42.681815
$ws.Range("O4").Value = 0.08068183420648613
$ws.Range("P4").Value = 0.08068183420648613
$ws.Range("Q4").Value = 23.88775985559334
$ws.Range("R4").Value = 214.98983870034
$ws.Range("S4").Value = 0.009505528070651571
$ws.Range("T4").Value = 0.009505528070651571

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Clcf1"
$ws.Range("C5").Value = "Il6st"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.679012
$ws.Range("H5").Value = 5.037036000000001
$ws.Range("I5").Value = 0.1178149724053671
$ws.Range("J5").Value = 0.1178149724053671
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 36.52018433333333
$ws.Range("N5").Value = 109.560553
$ws.Range("O5").Value = 0.20710333833547
$ws.Range("P5").Value = 0.2071033383354699
$ws.Range("Q5").Value = 61.31782773787867
$ws.Range("R5").Value = 551.8604496409081
$ws.Range("S5").Value = 0.02439987409105281
$ws.Range("T5").Value = 0.02439987409105281

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Clcf1"
$ws.Range("C6").Value = "Il6st"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.197979
$ws.Range("H6").Value = 9.593937
$ws.Range("I6").Value = 0.2243997110431275
$ws.Range("J6").Value = 0.2243997110431275
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 36.71344366666667
$ws.Range("N6").Value = 110.140331
$ws.Range("O6").Value = 0.2081992981130139
$ws.Range("P6").Value = 0.2081992981130138
$ws.Range("Q6").Value = 117.408821863683
$ws.Range("R6").Value = 1056.679396773147
$ws.Range("S6").Value = 0.04671986233594227
$ws.Range("T6").Value = 0.04671986233594225

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Clcf1"
$ws.Range("C7").Value = "Il6st"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.197979
$ws.Range("H7").Value = 9.593937
$ws.Range("I7").Value = 0.2243997110431275
$ws.Range("J7").Value = 0.2243997110431275
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 88.87708033333332
$ws.Range("N7").Value = 266.631241
$ws.Range("O7").Value = 0.5040155293450301
$ws.Range("P7").Value = 0.50401552934503
$ws.Range("Q7").Value = 284.227036487313
$ws.Range("R7").Value = 2558.043328385817
$ws.Range("S7").Value = 0.1131009391462737
$ws.Range("T7").Value = 0.1131009391462737

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Clcf1"
$ws.Range("C8").Value = "Il6st"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.197979
$ws.Range("H8").Value = 9.593937
$ws.Range("I8").Value = 0.2243997110431275
$ws.Range("J8").Value = 0.2243997110431275
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 14.22727166666667
$ws.Range("N8").Value = 42.681815
$ws.Range("O8").Value = 0.08068183420648613
$ws.Range("P8").Value = 0.08068183420648613
$ws.Range("Q8").Value = 45.498516017295
$ws.Range("R8").Value = 409.486644155655
$ws.Range("S8").Value = 0.01810498028236501
$ws.Range("T8").Value = 0.01810498028236501

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Clcf1"
$ws.Range("C9").Value = "Il6st"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.197979
$ws.Range("H9").Value = 9.593937
$ws.Range("I9").Value = 0.2243997110431275
$ws.Range("J9").Value = 0.2243997110431275
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 36.52018433333333
$ws.Range("N9").Value = 109.560553
$ws.Range("O9").Value = 0.20710333833547
$ws.Range("P9").Value = 0.2071033383354699
$ws.Range("Q9").Value = 116.790782574129
$ws.Range("R9").Value = 1051.117043167161
$ws.Range("S9").Value = 0.04647392927854653
$ws.Range("T9").Value = 0.04647392927854652

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Clcf1"
$ws.Range("C10").Value = "Il6st"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.919382666666667
$ws.Range("H10").Value = 5.758148
$ws.Range("I10").Value = 0.134681596027112
$ws.Range("J10").Value = 0.134681596027112
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 36.71344366666667
$ws.Range("N10").Value = 110.140331
$ws.Range("O10").Value = 0.2081992981130139
$ws.Range("P10").Value = 0.2081992981130138
$ws.Range("Q10").Value = 70.46714740744312
$ws.Range("R10").Value = 634.2043266669881
$ws.Range("S10").Value = 0.02804061376158519
$ws.Range("T10").Value = 0.02804061376158518

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Clcf1"
$ws.Range("C11").Value = "Il6st"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.919382666666667
$ws.Range("H11").Value = 5.758148
$ws.Range("I11").Value = 0.134681596027112
$ws.Range("J11").Value = 0.134681596027112
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 88.87708033333332
$ws.Range("N11").Value = 266.631241
$ws.Range("O11").Value = 0.5040155293450301
$ws.Range("P11").Value = 0.50401552934503
$ws.Range("Q11").Value = 170.5891274557409
$ws.Range("R11").Value = 1535.302147101668
$ws.Range("S11").Value = 0.06788161591463834
$ws.Range("T11").Value = 0.06788161591463833

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Clcf1"
$ws.Range("C12").Value = "Il6st"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.919382666666667
$ws.Range("H12").Value = 5.758148
$ws.Range("I12").Value = 0.134681596027112
$ws.Range("J12").Value = 0.134681596027112
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 14.22727166666667
$ws.Range("N12").Value = 42.681815
$ws.Range("O12").Value = 0.08068183420648613
$ws.Range("P12").Value = 0.08068183420648613
$ws.Range("Q12").Value = 27.30757863095778
$ws.Range("R12").Value = 245.76820767862
$ws.Range("S12").Value = 0.01086635820132439
$ws.Range("T12").Value = 0.01086635820132439

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Clcf1"
$ws.Range("C13").Value = "Il6st"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.919382666666667
$ws.Range("H13").Value = 5.758148
$ws.Range("I13").Value = 0.134681596027112
$ws.Range("J13").Value = 0.134681596027112
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 36.52018433333333
$ws.Range("N13").Value = 109.560553
$ws.Range("O13").Value = 0.20710333833547
$ws.Range("P13").Value = 0.2071033383354699
$ws.Range("Q13").Value = 70.09620879287155
$ws.Range("R13").Value = 630.865879135844
$ws.Range("S13").Value = 0.02789300814956406
$ws.Range("T13").Value = 0.02789300814956405

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Clcf1"
$ws.Range("C14").Value = "Il6st"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 7.454888
$ws.Range("H14").Value = 22.364664
$ws.Range("I14").Value = 0.5231037205243934
$ws.Range("J14").Value = 0.5231037205243932
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 36.71344366666667
$ws.Range("N14").Value = 110.140331
$ws.Range("O14").Value = 0.2081992981130139
$ws.Range("P14").Value = 0.2081992981130138
$ws.Range("Q14").Value = 273.6946106293094
$ws.Range("R14").Value = 2463.251495663784
$ws.Range("S14").Value = 0.1089098274534849
$ws.Range("T14").Value = 0.1089098274534848

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Clcf1"
$ws.Range("C15").Value = "Il6st"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 7.454888
$ws.Range("H15").Value = 22.364664
$ws.Range("I15").Value = 0.5231037205243934
$ws.Range("J15").Value = 0.5231037205243932
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 88.87708033333332
$ws.Range("N15").Value = 266.631241
$ws.Range("O15").Value = 0.5040155293450301
$ws.Range("P15").Value = 0.50401552934503
$ws.Range("Q15").Value = 662.5686796520026
$ws.Range("R15").Value = 5963.118116868024
$ws.Range("S15").Value = 0.2636523986024568
$ws.Range("T15").Value = 0.2636523986024567

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Clcf1"
$ws.Range("C16").Value = "Il6st"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 7.454888
$ws.Range("H16").Value = 22.364664
$ws.Range("I16").Value = 0.5231037205243934
$ws.Range("J16").Value = 0.5231037205243932
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 14.22727166666667
$ws.Range("N16").Value = 42.681815
$ws.Range("O16").Value = 0.08068183420648613
$ws.Range("P16").Value = 0.08068183420648613
$ws.Range("Q16").Value = 106.0627168205733
$ws.Range("R16").Value = 954.5644513851601
$ws.Range("S16").Value = 0.04220496765214516
$ws.Range("T16").Value = 0.04220496765214515

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Clcf1"
$ws.Range("C17").Value = "Il6st"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 7.454888
$ws.Range("H17").Value = 22.364664
$ws.Range("I17").Value = 0.5231037205243934
$ws.Range("J17").Value = 0.5231037205243932
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 36.52018433333333
$ws.Range("N17").Value = 109.560553
$ws.Range("O17").Value = 0.20710333833547
$ws.Range("P17").Value = 0.2071033383354699
$ws.Range("Q17").Value = 272.2538839443547
$ws.Range("R17").Value = 2450.284955499192
$ws.Range("S17").Value = 0.1083365268163066
$ws.Range("T17").Value = 0.1083365268163065
